$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Channel Assignments")
$ws.Range("C42").Value = "mastLeft"
$ws.Range("E42").Value = "talon"
$ws.Range("C43").Value = "mastRight"
$ws.Range("E43").Value = "talon"
$ws.Range("C65").Value = "mastPot"
